$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 18:35"

# Estados Unidos (row 4) - updated counts
$ws.Range("B4").Value = 1414053
$ws.Range("C4").Value = 5417
$ws.Range("E4").Value = 1031276
$ws.Range("G4").Value = 344
$ws.Range("H4").Value = 83769

# Italia (row 8) - updated counts
$ws.Range("B8").Value = 222104
$ws.Range("C8").Value = 888
$ws.Range("D8").Value = 112541
$ws.Range("E8").Value = 78457
$ws.Range("F8").Value = 893
$ws.Range("G8").Value = 195
$ws.Range("H8").Value = 31106

# Turquia (row 12) - updated counts
$ws.Range("B12").Value = 143114
$ws.Range("C12").Value = 1639
$ws.Range("D12").Value = 101715
$ws.Range("E12").Value = 37447
$ws.Range("F12").Value = 998
$ws.Range("G12").Value = 58
$ws.Range("H12").Value = 3952

# Egipto overtakes Serbia (rows 48-49 swap with Egipto's updated counts)
$ws.Range("A48").Value = "Egipto"
$ws.Range("B48").Value = 10431
$ws.Range("C48").Value = 338
$ws.Range("D48").Value = 2486
$ws.Range("E48").Value = 7389
$ws.Range("F48").Value = 41
$ws.Range("G48").Value = 12
$ws.Range("H48").Value = 556

$ws.Range("A49").Value = "Serbia"
$ws.Range("B49").Value = 10295
$ws.Range("C49").Value = 52
$ws.Range("D49").Value = 3824
$ws.Range("E49").Value = 6249
$ws.Range("F49").Value = 22
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 222

# Suazilandia overtakes Martinica and Islas Feroe (rows 149-151 reorder with Suazilandia's updated counts)
$ws.Range("A149").Value = "Suazilandia"
$ws.Range("B149").Value = 187
$ws.Range("C149").Value = 3
$ws.Range("D149").Value = 48
$ws.Range("E149").Value = 137
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 2

$ws.Range("A150").Value = "Martinica"
$ws.Range("B150").Value = 187
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 91
$ws.Range("E150").Value = 82
$ws.Range("F150").Value = 4
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 14

$ws.Range("A151").Value = "Islas Feroe"
$ws.Range("B151").Value = 187
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 187
$ws.Range("E151").Value = 0
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0

# Nueva Caledonia and Belice swap order (rows 193-194, no value changes otherwise)
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

$ws.Range("A194").Value = "Belice"
$ws.Range("B194").Value = 18
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 16
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2
